$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 1.7
$ws.Range("F4").Value = 99

$ws.Range("E5").Value = 1.7
$ws.Range("F5").Value = 98

$ws.Range("F6").Value = 99

$ws.Range("F7").Value = 99

$ws.Range("F8").Value = 98

$ws.Range("E6").Select()
